{"js": "// Change 1: \"Model AIC: 0.18\" -> \"Model AIC: 487.69\"\nconst aicModelHits = context.document.body.search(\"Model AIC: 0.18\", { matchCase: true });\naicModelHits.load(\"items\");\nawait context.sync();\nif (aicModelHits.items.length > 0) {\n  aicModelHits.items[0].insertText(\"Model AIC: 487.69\", Word.InsertLocation.replace);\n}\n\n// Change 2: \" crude odds ratio\" -> \" crude odds ratio, \" + italic \"AIC\" + \" Akaike Information Criterion\"\nconst crudeHits = context.document.body.search(\" crude odds ratio\", { matchCase: true });\ncrudeHits.load(\"items\");\nawait context.sync();\nif (crudeHits.items.length > 0) {\n  const target = crudeHits.items[0];\n  const newRange = target.insertText(\n    \" crude odds ratio, AIC Akaike Information Criterion\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n\n  const aicHits = newRange.search(\"AIC\", { matchCase: true });\n  aicHits.load(\"items\");\n  await context.sync();\n  if (aicHits.items.length > 0) {\n    aicHits.items[0].font.italic = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: \"Model AIC: 0.18\" -> \"Model AIC: 487.69\"\n$find1 = $d.Content.Find\n$found1 = $find1.Execute(\"Model AIC: 0.18\", $false, $false, $false, $false, $false, $true, 1, $false, \"Model AIC: 487.69\", 2)\n\n# Change 2: \" crude odds ratio\" -> \" crude odds ratio, \" + italic \"AIC\" + \" Akaike Information Criterion\"\n$find2 = $d.Content.Find\n$find2.Text = \" crude odds ratio\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $r = $find2.Parent\n    $origStart = $r.Start\n    $replacement = \" crude odds ratio, AIC Akaike Information Criterion\"\n    $r.Text = $replacement\n\n    # Compute the exact character span of \"AIC\" within the freshly written\n    # text (by offset, since Range.Text reads back stale after a .Text set)\n    # and italicize only that span, leaving the rest of the run untouched.\n    $aicOffset = $replacement.IndexOf(\"AIC\")\n    $aicStart = $origStart + $aicOffset\n    $aicEnd = $aicStart + 3\n    $aicRange = $d.Range($aicStart, $aicEnd)\n    $aicRange.Font.Italic = 1\n}\n\nWrite-Output \"done\"\n"}
